# AutoCommit_11 апреля 2024 г. 15:58:26_SibNout2023
#
# - Highlight the per-homework score cells (columns C:E, rows 4-32) with a
#   solid green fill, mirroring borderId=1 already in place.
# - Normalize previously-blank score cells in C:E to an explicit 0.
# - Add a new (unlabeled) column N, one row below M, filled with 0 for every
#   student row.
# - Restore the "working" view: 145% zoom, frozen pane scrolled back up to
#   C4, and the active selection on H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow  = 32

# Excel's RGB(0,176,80) == the "FF00B050" green used for the fill.
$green = 5287936

# Columns C (3) through E (5) that used to render as blank (self-closed <c/>)
# and should now hold an explicit 0. Rows not listed here already had a
# numeric value in every one of C/D/E and are left untouched.
$emptyCols = @{
    4  = @("C", "D", "E")
    5  = @("C", "D", "E")
    9  = @("C", "D", "E")
    12 = @("C", "D", "E")
    13 = @("C", "D", "E")
    16 = @("C", "D", "E")
    17 = @("C", "E")
    20 = @("C", "D", "E")
    22 = @("E")
    23 = @("C", "D", "E")
    26 = @("C", "D", "E")
    27 = @("C", "D", "E")
    28 = @("C", "D", "E")
    30 = @("D", "E")
    32 = @("C", "D", "E")
}

for ($r = $firstRow; $r -le $lastRow; $r++) {

    # C and D always pick up the new green fill (style -> fontId0/fillId2/
    # borderId1). E only does for rows 4-29 -- rows 30-32 keep their old
    # (unfilled) border-only style, matching the source edit exactly.
    $ws.Range("C$r`:D$r").Interior.Color = $green
    if ($r -le 29) {
        $ws.Range("E$r").Interior.Color = $green
    }

    # Fill in explicit zeros for cells that used to be blank.
    if ($emptyCols.ContainsKey($r)) {
        foreach ($col in $emptyCols[$r]) {
            $ws.Range("$col$r").Value = 0
        }
    }

    # New column N: zero for every student row.
    $ws.Range("N$r").Value = 0
}

# View: 145% zoom, scrolled back to the top of the frozen area, selection on H8.
$excel.ActiveWindow.Zoom = 145
$ws.Range("H8").Select()
